# "Excel com tempo de duração dos testes"
#
# Adds a "runtime" column (G) and a merged-looking title/banner row above the
# existing header row, explaining that the attack duration (in seconds) will
# always be greater than the sum of each prompt's "runtime".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row above row 1. This pushes the existing header
#    row (and every row below it) down by one, exactly like pressing
#    "Insert Sheet Rows" with row 1 selected in Excel.
$ws.Rows.Item(1).EntireRow.Insert()

# 2) Insert a brand-new column at G for the "runtime" data the team now
#    wants to capture per prompt.
$ws.Columns.Item(7).EntireColumn.Insert()

# 3) Header label for the new column, in the (now shifted-down) header row.
$ws.Cells.Item(2, 7).Value = "runtime"

# 4) Build the new banner/title cell in A1 explaining the new column and
#    warning that total duration > sum of per-prompt runtimes.
$titleCell = $ws.Cells.Item(1, 1)
$titleCell.Value = "Duração do ataque, em segundos: " + [char]10 + "(Atenção! A duração do ataque vai ser sempre superior à soma do 'runtime' de cada prompt)"

# Bold font, yellow fill, centered + wrapped text - a highlighted banner.
$titleCell.Font.Bold = $true
$titleCell.Interior.Color = 65535    # RGB(255,255,0) -> yellow
$titleCell.HorizontalAlignment = -4108  # xlCenter
$titleCell.VerticalAlignment = -4108    # xlCenter
$titleCell.WrapText = $true

# Give the banner row enough height to show the two wrapped lines of text.
$ws.Rows.Item(1).RowHeight = 37.2

# 5) Widen columns A and B a bit (and size the new runtime column G) so the
#    new/longer text fits comfortably.
$ws.Columns.Item(1).ColumnWidth = 77.8
$ws.Columns.Item(2).ColumnWidth = 25.97
$ws.Columns.Item(7).ColumnWidth = 16.4

# 6) Leave the cursor where the author left it when saving.
$ws.Range("A7").Select()
